# The diff shows a new price record being inserted into the "Zapallo" sheet
# at row 61 (pushing the existing rows 61-141 down to 62-142, and growing the
# used range from A1:R141 to A1:R142). Reproduce that with a real row insert
# so every downstream row's data shifts down exactly one position, then
# populate the newly-opened row 61 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 61; this pushes old rows 61..141 down to 62..142.
$ws.Rows.Item(61).Insert()

# Fill the new row 61 with the inserted record.
$ws.Range("A61").Value = 7
$ws.Range("B61").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C61").Value = 'Ñuble'
$ws.Range("D61").Value = 44650
$ws.Range("E61").Value = 16
$ws.Range("F61").Value = 100112045
$ws.Range("G61").Value = 'Zapallo'
$ws.Range("H61").Value = 'Camote'
$ws.Range("I61").Value = '1a (cosecha)'
$ws.Range("J61").Value = 200
$ws.Range("K61").Value = 300
$ws.Range("L61").Value = 350
$ws.Range("M61").Value = 325
$ws.Range("N61").Value = '$/kilo (volumen en unidades)'
$ws.Range("O61").Value = "Región de O'Higgins"
$ws.Range("P61").Value = 325
$ws.Range("Q61").Value = 1
$ws.Range("R61").Value = 'Hortaliza'
